$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.062.56"
$ws.Range("E2").Value = "  +1.70%  "
# Row 3
$ws.Range("D3").Value = "3.339.52"
$ws.Range("E3").Value = "  +1.90%  "
# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").Value = "'582.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.14%  "
# Row 6
$ws.Range("D6").Value = "'177.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.00%  "
# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "
# Row 8
$ws.Range("D8").Value = "'0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.62%  "
# Row 9
$ws.Range("D9").Value = "3.337.02"
$ws.Range("E9").Value = "  +2.13%  "
# Row 10
$ws.Range("D10").Value = "'0.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.45%  "
# Row 11
$ws.Range("D11").Value = "'0.582"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.21%  "
# Row 12
$ws.Range("D12").Value = "'47.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.99%  "
# Row 13
$ws.Range("D13").Value = "'0.0000274"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.54%  "
# Row 14
$ws.Range("D14").Value = "'686.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.06%  "
# Row 15
$ws.Range("D15").Value = "3.884.55"
$ws.Range("E15").Value = "  +2.07%  "
# Row 16
$ws.Range("D16").Value = "'8.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.53%  "
# Row 17
$ws.Range("D17").Value = "68.095.97"
$ws.Range("E17").Value = "  +1.53%  "
# Row 18
$ws.Range("D18").Value = "'0.118"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
# Row 19
$ws.Range("D19").Value = "3.341.96"
$ws.Range("E19").Value = "  +1.73%  "
# Row 20
$ws.Range("D20").Value = "'17.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.37%  "
# Row 21
$ws.Range("D21").Value = "'11.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.76%  "
# Row 22
$ws.Range("D22").Value = "'0.900"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.03%  "
# Row 23
$ws.Range("D23").Value = "'5.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.19%  "
# Row 24
$ws.Range("D24").Value = "'17.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.48%  "
# Row 25
$ws.Range("D25").Value = "'99.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "
# Row 26
$ws.Range("D26").Value = "'3.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.19%  "
# Row 27
$ws.Range("D27").Value = "'2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.06%  "
# Row 28
$ws.Range("D28").Value = "'9.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.56%  "
# Row 29
$ws.Range("D29").Value = "'33.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.39%  "
# Row 30
$ws.Range("D30").Value = "'8.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.30%  "
# Row 31
$ws.Range("D31").Value = "'7.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.15%  "
# Row 32
$ws.Range("D32").Value = "'571.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.08%  "
# Row 33
$ws.Range("D33").Value = "'11.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.49%  "
# Row 34
$ws.Range("D34").Value = "'0.106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.19%  "
# Row 35
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
# Row 36
$ws.Range("D36").Value = "'57.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.64%  "
# Row 37
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.717.56"
$ws.Range("E37").Value = "  -3.95%  "
# Row 38
$ws.Range("D38").Value = "'3.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.58%  "
# Row 39
$ws.Range("D39").Value = "'34.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.88%  "
# Row 40
$ws.Range("D40").Value = "'0.132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.84%  "
# Row 41
$ws.Range("D41").Value = "'3.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.82%  "
# Row 42
$ws.Range("D42").Value = "'2.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.26%  "
# Row 43
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "'3.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "
# Row 44
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0678"
$ws.Range("E44").Value = "  +1.93%  "
# Row 45
$ws.Range("D45").Value = "'0.338"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.80%  "
# Row 46
$ws.Range("D46").Value = "'0.0408"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.14%  "
# Row 47
$ws.Range("D47").Value = "'2.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.81%  "
# Row 48
$ws.Range("D48").Value = "'0.129"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "
# Row 49
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.37%  "
# Row 50
$ws.Range("D50").Value = "'1.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.07%  "
# Row 51
$ws.Range("D51").Value = "'129.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
